$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 6;  C = 9; D = 0; E = 6 }
    3  = @{ B = 7;  C = 9; D = 0; E = 7 }
    4  = @{ B = 5;  C = 9; D = 0; E = 5 }
    6  = @{ B = 16; C = 1; D = 8; E = 0 }
    7  = @{ B = 6;  C = 9; D = 0; E = 6 }
    8  = @{ B = 2;  C = 9; D = 0; E = 2 }
    9  = @{ B = 18; C = 0; D = 9; E = 0 }
    10 = @{ B = 16; C = 1; D = 8; E = 0 }
    11 = @{ B = 8;  C = 9; D = 0; E = 8 }
    12 = @{ B = 6;  C = 9; D = 0; E = 6 }
    13 = @{ B = 6;  C = 9; D = 0; E = 6 }
    14 = @{ B = 5;  C = 9; D = 0; E = 5 }
    15 = @{ B = 6;  C = 9; D = 0; E = 6 }
    16 = @{ B = 1;  C = 9; D = 0; E = 1 }
    17 = @{ B = 8;  C = 9; D = 0; E = 8 }
    18 = @{ B = 4;  C = 9; D = 0; E = 4 }
    19 = @{ B = 16; C = 1; D = 8; E = 0 }
    20 = @{ B = 16; C = 1; D = 8; E = 0 }
    21 = @{ B = 14; C = 0; D = 7; E = 0 }
    22 = @{ B = 2;  C = 9; D = 0; E = 2 }
    23 = @{ B = 12; C = 3; D = 6; E = 0 }
}

foreach ($rowNum in $data.Keys) {
    $row = $data[$rowNum]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$rowNum").Value = $row[$col]
    }
}
